$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 72380.14
$ws.Range("I11").Value = 72380.14
$ws.Range("K11").Value = 72380.14
$ws.Range("M11").Value = -72240.14
$ws.Range("H70").Value = 4385.654
$ws.Range("I70").Value = 1189.2307
$ws.Range("J70").Value = 7582.077
$ws.Range("K70").Value = 3567.6921
$ws.Range("L70").Value = 22746.231
$ws.Range("M70").Value = -3297.6921
$ws.Range("N70").Value = -23286.231
$ws.Range("H73").Value = 4385.654
$ws.Range("I73").Value = 1189.2307
$ws.Range("J73").Value = 7582.077
$ws.Range("K73").Value = 3567.6921
$ws.Range("L73").Value = 22746.231
$ws.Range("M73").Value = -2631.6921
$ws.Range("N73").Value = -24618.231
$ws.Range("H137").Value = 22059790
$ws.Range("I137").Value = 5953333.5
$ws.Range("J137").Value = 55883348
$ws.Range("K137").Value = 17860000.5
$ws.Range("L137").Value = 167650044
$ws.Range("M137").Value = -17857450.5
$ws.Range("N137").Value = -167655144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 44976548
$ws.Range("I74").Value = 78023290
$ws.Range("J74").Value = 14290285
$ws.Range("K74").Value = 78023290
$ws.Range("L74").Value = 14290285
$ws.Range("M74").Value = -78022416
$ws.Range("N74").Value = -14292033
$ws.Range("H77").Value = 44976548
$ws.Range("I77").Value = 78023290
$ws.Range("J77").Value = 14290285
$ws.Range("K77").Value = 390116450
$ws.Range("L77").Value = 71451425
$ws.Range("M77").Value = -390112082
$ws.Range("N77").Value = -71460161
$ws.Range("H88").Value = 4729.8667
$ws.Range("I88").Value = 2175
$ws.Range("J88").Value = 5658.909
$ws.Range("K88").Value = 2175
$ws.Range("L88").Value = 5658.909
$ws.Range("M88").Value = -1769
$ws.Range("N88").Value = -6470.909
$ws.Range("H91").Value = 4729.8667
$ws.Range("I91").Value = 2175
$ws.Range("J91").Value = 5658.909
$ws.Range("K91").Value = 2175
$ws.Range("L91").Value = 5658.909
$ws.Range("M91").Value = -771
$ws.Range("N91").Value = -8466.909
$ws.Range("H122").Value = 2132.2727
$ws.Range("I122").Value = 1601.5385
$ws.Range("J122").Value = 2898.889
$ws.Range("K122").Value = 4804.6155
$ws.Range("L122").Value = 8696.667000000001
$ws.Range("M122").Value = -2354.6155
$ws.Range("N122").Value = -13596.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1842.4362
$ws.Range("I86").Value = 1886.236
$ws.Range("J86").Value = 1062.8
$ws.Range("K86").Value = 1886.236
$ws.Range("L86").Value = 1062.8
$ws.Range("M86").Value = -763.2360000000001
$ws.Range("N86").Value = -3308.8
$ws.Range("H89").Value = 1842.4362
$ws.Range("I89").Value = 1886.236
$ws.Range("J89").Value = 1062.8
$ws.Range("K89").Value = 9431.18
$ws.Range("L89").Value = 5314
$ws.Range("M89").Value = -3815.18
$ws.Range("N89").Value = -16546
$ws.Range("H134").Value = 35715640
$ws.Range("I134").Value = 38462984
$ws.Range("J134").Value = 17857892
$ws.Range("K134").Value = 115388952
$ws.Range("L134").Value = 53573676
$ws.Range("M134").Value = -115386417
$ws.Range("N134").Value = -53578746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3171.4285
$ws.Range("I62").Value = 2618.182
$ws.Range("K62").Value = 2618.182
$ws.Range("M62").Value = -1994.182
$ws.Range("H65").Value = 3171.4285
$ws.Range("I65").Value = 2618.182
$ws.Range("K65").Value = 13090.91
$ws.Range("M65").Value = -9970.91
$ws.Range("H132").Value = 2501675.2
$ws.Range("I132").Value = 3572428
$ws.Range("K132").Value = 10717284
$ws.Range("M132").Value = -10714754

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2903
$ws.Range("J117").Value = 3589.5715
$ws.Range("L117").Value = 10768.7145
$ws.Range("N117").Value = -17652.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9262398
$ws.Range("I122").Value = 3306
$ws.Range("K122").Value = 9918
$ws.Range("M122").Value = -7468
$ws.Range("H132").Value = 13022918
$ws.Range("I132").Value = 19049238
$ws.Range("J132").Value = 6996598.5
$ws.Range("K132").Value = 57147714
$ws.Range("L132").Value = 20989795.5
$ws.Range("M132").Value = -57145184
$ws.Range("N132").Value = -20994855.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 91865.73
$ws.Range("I16").Value = 125440.375
$ws.Range("J16").Value = 2333.3333
$ws.Range("K16").Value = 125440.375
$ws.Range("L16").Value = 2333.3333
$ws.Range("M16").Value = -125270.375
$ws.Range("N16").Value = -2673.3333
$ws.Range("H22").Value = 5578
$ws.Range("I22").Value = 447.5
$ws.Range("J22").Value = 9682.4
$ws.Range("K22").Value = 447.5
$ws.Range("L22").Value = 9682.4
$ws.Range("M22").Value = -152.5
$ws.Range("N22").Value = -10272.4
$ws.Range("H27").Value = 5578
$ws.Range("I27").Value = 447.5
$ws.Range("J27").Value = 9682.4
$ws.Range("K27").Value = 447.5
$ws.Range("L27").Value = 9682.4
$ws.Range("M27").Value = -340.5
$ws.Range("N27").Value = -9896.4
$ws.Range("H82").Value = 6033.4546
$ws.Range("I82").Value = 2152.75
$ws.Range("J82").Value = 8251
$ws.Range("K82").Value = 2152.75
$ws.Range("L82").Value = 8251
$ws.Range("M82").Value = -1791.75
$ws.Range("N82").Value = -8973
$ws.Range("H85").Value = 6033.4546
$ws.Range("I85").Value = 2152.75
$ws.Range("J85").Value = 8251
$ws.Range("K85").Value = 2152.75
$ws.Range("L85").Value = 8251
$ws.Range("M85").Value = -904.75
$ws.Range("N85").Value = -10747
$ws.Range("H122").Value = 8513543
$ws.Range("I122").Value = 1066435.2
$ws.Range("J122").Value = 33337234
$ws.Range("K122").Value = 3199305.6
$ws.Range("L122").Value = 100011702
$ws.Range("M122").Value = -3196855.6
$ws.Range("N122").Value = -100016602
$ws.Range("H132").Value = 2566450.2
$ws.Range("I132").Value = 3510784.2
$ws.Range("J132").Value = 3257.5715
$ws.Range("K132").Value = 10532352.6
$ws.Range("L132").Value = 9772.7145
$ws.Range("M132").Value = -10529822.6
$ws.Range("N132").Value = -14832.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 14158.292
$ws.Range("I81").Value = 2075.25
$ws.Range("J81").Value = 16574.9
$ws.Range("K81").Value = 4150.5
$ws.Range("L81").Value = 33149.8
$ws.Range("M81").Value = -3089.5
$ws.Range("N81").Value = -35271.8
$ws.Range("H84").Value = 14158.292
$ws.Range("I84").Value = 2075.25
$ws.Range("J84").Value = 16574.9
$ws.Range("K84").Value = 20752.5
$ws.Range("L84").Value = 165749
$ws.Range("M84").Value = -15448.5
$ws.Range("N84").Value = -176357
$ws.Range("H122").Value = 1538.0714
$ws.Range("I122").Value = 1388.8334
$ws.Range("K122").Value = 4166.5002
$ws.Range("M122").Value = -1716.5002
$ws.Range("H132").Value = 2027442
$ws.Range("I132").Value = 1645435
$ws.Range("J132").Value = 2756728
$ws.Range("K132").Value = 4936305
$ws.Range("L132").Value = 8270184
$ws.Range("M132").Value = -4933775
$ws.Range("N132").Value = -8275244
$ws.Range("H136").Value = 2402.8696
$ws.Range("I136").Value = 994.55554
$ws.Range("J136").Value = 3308.2144
$ws.Range("K136").Value = 2983.66662
$ws.Range("L136").Value = 9924.643199999999
$ws.Range("M136").Value = -433.66662
$ws.Range("N136").Value = -15024.6432
